$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $value) {
    $r = $ws.Range($cellAddr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "65.705.78"
$ws.Range("E2").Value = "  -2.52%  "
Set-TextValue "D3" "3.282.29"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue "D5" "573.06"
$ws.Range("E5").Value = "  -0.84%  "
Set-TextValue "D6" "177.57"
$ws.Range("E6").Value = "  -4.47%  "
Set-TextValue "D7" "0.632"
$ws.Range("E7").Value = "  +4.72%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -2.65%  "
Set-TextValue "D10" "6.71"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E11").Value = "  -2.74%  "
Set-TextValue "D12" "3.854.44"
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("E13").Value = "  -3.63%  "
Set-TextValue "D14" "26.60"
$ws.Range("E14").Value = "  -3.14%  "
Set-TextValue "D15" "65.807.86"
$ws.Range("E15").Value = "  -2.70%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D16" "0.0000163"
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "3.293.07"
$ws.Range("E17").Value = "  -0.63%  "
Set-TextValue "D18" "437.03"
$ws.Range("E18").Value = "  -1.42%  "
Set-TextValue "D19" "5.58"
$ws.Range("E19").Value = "  -2.52%  "
Set-TextValue "D20" "13.22"
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("E21").Value = "  -4.68%  "
Set-TextValue "D22" "72.47"
$ws.Range("E22").Value = "  -2.05%  "
Set-TextValue "D23" "1.00"
$ws.Range("E23").Value = "  +0.25%  "
Set-TextValue "D24" "3.432.21"
$ws.Range("E24").Value = "  -0.69%  "
Set-TextValue "D25" "0.510"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("E26").Value = "  -5.00%  "
$ws.Range("E27").Value = "  +3.03%  "
Set-TextValue "D28" "8.91"
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -2.12%  "
Set-TextValue "D31" "22.33"
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("E32").Value = "  +0.12%  "
Set-TextValue "D33" "5.15"
$ws.Range("E33").Value = "  -3.67%  "
Set-TextValue "D34" "6.63"
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("E35").Value = "  -4.07%  "
Set-TextValue "D36" "158.64"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("E37").Value = "  -4.89%  "
Set-TextValue "D38" "26.82"
$ws.Range("E38").Value = "  -1.35%  "
Set-TextValue "D39" "1.79"
$ws.Range("E39").Value = "  -3.63%  "
Set-TextValue "D40" "2.775.77"
$ws.Range("E40").Value = "  +0.25%  "
Set-TextValue "D41" "0.779"
$ws.Range("E41").Value = "  -1.50%  "
Set-TextValue "D42" "4.33"
$ws.Range("E42").Value = "  -3.42%  "
Set-TextValue "D43" "40.31"
$ws.Range("E43").Value = "  +0.37%  "
Set-TextValue "D44" "6.05"
$ws.Range("E44").Value = "  -3.20%  "
$ws.Range("E45").Value = "  -2.43%  "
Set-TextValue "D46" "2.29"
$ws.Range("E46").Value = "  -4.87%  "
Set-TextValue "D47" "319.64"
$ws.Range("E47").Value = "  -2.64%  "
Set-TextValue "D48" "23.41"
$ws.Range("E48").Value = "  -5.80%  "
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("E50").Value = "  +2.27%  "
Set-TextValue "D51" "0.999"
$ws.Range("E51").Value = "  -0.04%  "
